$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply scraped crypto market data refresh (prices & 1h volume deltas).
# Values are forced to text via a leading apostrophe (mirrors the source
# data which stores prices/percentages as literal strings, e.g. "0.07100"
# or "26.321.49"), then the style is reset to Normal so no stray number
# format / quote-prefix styling is left behind on the cell.

$ws.Range('D2').Value = "'26.321.49"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -2.03%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.791.74"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -2.12%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'1.008"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.09%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'1.007"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +0.03%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'307.51"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -0.96%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.4526"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  -2.11%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.3591"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  -2.67%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'45.59"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -0.47%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.07100"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -1.13%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.8844"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +0.88%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'0.07819"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -0.39%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'19.48"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -0.57%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'1.814.33"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -2.56%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'5.285"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -0.92%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'6.327"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -0.97%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'84.63"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -2.56%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('E18').Value = "'  +0.22%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'0.000008538"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E20').Value = "'  +0.06%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'14.25"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -1.52%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'26.361.37"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -1.97%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'4.983"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -0.17%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('B24').Value = "'WrappedliquidstakedEther2.0"
$ws.Range('B24').Style = 'Normal'
$ws.Range('C24').Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range('C24').Style = 'Normal'
$ws.Range('D24').Value = "'2.023.13"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -2.31%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('B25').Value = "'Cosmos"
$ws.Range('B25').Style = 'Normal'
$ws.Range('C25').Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range('C25').Style = 'Normal'
$ws.Range('D25').Value = "'10.49"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +0.56%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'1.975"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -0.14%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'152.17"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +0.81%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'17.87"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -1.94%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'2.027"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  +3.31%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'111.92"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -1.50%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'4.850"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -1.57%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'0.08670"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  -1.67%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'3.061"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -2.22%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'4.445"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -0.38%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'0.7240"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -3.85%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'2.703"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +5.88%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'1.105"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -2.29%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'1.072"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -1.46%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'0.01931"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -0.16%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = "'  -0.55%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'2.871"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -2.06%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'0.5093"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +2.38%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'6.853"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -0.69%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'0.1514"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -5.17%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'7.982"
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Value = "'1.007"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -0.01%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'0.4629"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -1.05%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('B48').Value = "'Quant"
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').Value = "'100.73"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -1.51%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').Value = "'EnergySwap"
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').Value = "'9.831"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -2.69%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'1.579"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -1.84%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'0.05969"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -2.12%  "
$ws.Range('E51').Style = 'Normal'
